$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Renumber the "SlNo" column (A) for rows 5-19 so the sequence is
# contiguous (4..18) instead of having gaps (5,6,8,9,10,12,13,14,16,17,18,20,17,18,19).
$ws.Range("A5").Value = 4
$ws.Range("A6").Value = 5
$ws.Range("A7").Value = 6
$ws.Range("A8").Value = 7
$ws.Range("A9").Value = 8
$ws.Range("A10").Value = 9
$ws.Range("A11").Value = 10
$ws.Range("A12").Value = 11
$ws.Range("A13").Value = 12
$ws.Range("A14").Value = 13
$ws.Range("A15").Value = 14
$ws.Range("A16").Value = 15
$ws.Range("A17").Value = 16
$ws.Range("A18").Value = 17
$ws.Range("A19").Value = 18

# Update the sheet's saved selection to span the whole renumbered column.
$ws.Range("A2:A19").Select()
